# Update the "want-to-go" counts (column F) and the minimum ticket price for
# row 2 (column G) on the "展览" and "全部类型" sheets, reflecting the latest
# generated data snapshot (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F ("想去人数")
$fUpdates = @{
    5  = 23
    6  = 115
    7  = 90
    9  = 49
    11 = 576
    13 = 301
    15 = 373
    17 = 92
    18 = 9
    19 = 53
    20 = 49
    21 = 100
    22 = 944
    23 = 1402
    25 = 328
    28 = 156
    29 = 42
    31 = 221
    32 = 252
    34 = 1623
    36 = 99
    38 = 584
    40 = 3682
    41 = 427
    42 = 203
    46 = 67
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Column G row 2: minimum ticket price 45 -> 55
    $ws.Range("G2").Value = 55

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }
}
